$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A4:A5").NumberFormat = "@"

$ws.Cells.Item(4, 1).Value = "2025-09-13"
$ws.Cells.Item(4, 2).Value = "AAA"
$ws.Cells.Item(4, 3).Value = "44CDX012"
$ws.Cells.Item(4, 4).Value = "MAMA Nagar"

$ws.Cells.Item(5, 1).Value = "2025-08-13"
$ws.Cells.Item(5, 2).Value = "JJJ"
$ws.Cells.Item(5, 3).Value = "456CDX0176"
$ws.Cells.Item(5, 4).Value = "AMMA Nagar"

$ws.Range("A4:A5").ClearFormats()
